$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Usage" values between the BTN pins (tested configuration update)
$h14 = $ws.Range("H14").Value()
$h15 = $ws.Range("H15").Value()
$h17 = $ws.Range("H17").Value()
$h18 = $ws.Range("H18").Value()

$ws.Range("H14").Value = $h18
$ws.Range("H15").Value = $h17
$ws.Range("H17").Value = $h15
$ws.Range("H18").Value = $h14

# Update the active cell selection to match the saved view
$ws.Range("H9").Select()
